# Update cryptocurrency price/volume data to reflect the latest scrape.
# (Two coin pairs - Avalanche/TRON and Fetch.AI/USDe - also swapped rank position.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.546.51'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').Value = '3.387.36'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''575.34'
$ws.Range('E5').Value = '  +0.36%  '
$ws.Range('D6').Value = '''140.53'
$ws.Range('E6').Value = '  -1.60%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -0.53%  '
$ws.Range('D9').Value = '''7.71'
$ws.Range('E9').Value = '  +2.00%  '
$ws.Range('D10').Value = '''0.122'
$ws.Range('E10').Value = '  -1.24%  '
$ws.Range('E11').Value = '  -2.52%  '
$ws.Range('D12').Value = '3.969.04'
$ws.Range('E12').Value = '  -0.18%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '''0.125'
$ws.Range('E13').Value = '  +0.12%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').Value = '''28.44'
$ws.Range('E14').Value = '  +0.91%  '
$ws.Range('D15').Value = '3.403.72'
$ws.Range('E15').Value = '  +0.34%  '
$ws.Range('E16').Value = '  -0.70%  '
$ws.Range('D17').Value = '61.554.87'
$ws.Range('E17').Value = '  +0.80%  '
$ws.Range('D18').Value = '''6.12'
$ws.Range('E18').Value = '  -0.50%  '
$ws.Range('D19').Value = '''13.62'
$ws.Range('E19').Value = '  -1.69%  '
$ws.Range('D20').Value = '''8.99'
$ws.Range('E20').Value = '  +0.14%  '
$ws.Range('D21').Value = '''391.17'
$ws.Range('E21').Value = '  +2.05%  '
$ws.Range('D22').Value = '''75.23'
$ws.Range('E22').Value = '  +1.32%  '
$ws.Range('D23').Value = '''0.553'
$ws.Range('E23').Value = '  -1.24%  '
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('E25').Value = '  -5.23%  '
$ws.Range('D26').Value = '''0.193'
$ws.Range('E26').Value = '  +7.13%  '
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('D28').Value = '''7.24'
$ws.Range('E28').Value = '  -2.13%  '
$ws.Range('D29').Value = '''8.04'
$ws.Range('E29').Value = '  +0.31%  '
$ws.Range('D30').Value = '''2.14'
$ws.Range('E30').Value = '  -0.77%  '
$ws.Range('B31').Value = 'USDe'
$ws.Range('C31').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D31').Value = '''1.00'
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').Value = '''1.39'
$ws.Range('E32').Value = '  -2.07%  '
$ws.Range('D33').Value = '''23.26'
$ws.Range('E33').Value = '  -1.11%  '
$ws.Range('D34').Value = '''6.91'
$ws.Range('E34').Value = '  -1.42%  '
$ws.Range('D35').Value = '''168.35'
$ws.Range('E35').Value = '  +0.30%  '
$ws.Range('D36').Value = '''5.04'
$ws.Range('E36').Value = '  +0.69%  '
$ws.Range('D37').Value = '3.422.81'
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('E38').Value = '  -1.69%  '
$ws.Range('D39').Value = '''0.0767'
$ws.Range('E39').Value = '  -0.71%  '
$ws.Range('D40').Value = '''26.11'
$ws.Range('E40').Value = '  -6.21%  '
$ws.Range('D41').Value = '''0.778'
$ws.Range('E41').Value = '  -0.36%  '
$ws.Range('E42').Value = '  -0.60%  '
$ws.Range('E43').Value = '  -1.38%  '
$ws.Range('E44').Value = '  +2.32%  '
$ws.Range('D45').Value = '2.454.30'
$ws.Range('E45').Value = '  -0.99%  '
$ws.Range('D46').Value = '''22.94'
$ws.Range('E46').Value = '  -0.28%  '
$ws.Range('D47').Value = '''6.65'
$ws.Range('E47').Value = '  -2.54%  '
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('D50').Value = '''2.02'
$ws.Range('E50').Value = '  -4.82%  '
$ws.Range('E51').Value = '  -2.13%  '
